$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ENVELOPE_HEADER_TALLYREQUEST"
$ws.Range("B1").Value = "ENVELOPE_BODY_IMPORTDATA_REQUESTDESC_REPORTNAME"
$ws.Range("C1").Value = "ENVELOPE_BODY_IMPORTDATA_REQUESTDESC_STATICVARIABLES_SVCURRENTCOMPANY"
$ws.Range("D1").Value = "ENVELOPE_BODY_IMPORTDATA_REQUESTDATA_TALLYMESSAGE"
